# The workbook tracks weekly price data for "Vega Monumental Concepción - Zapallo".
# This edit inserts one new weekly record at row 198 (pushing the existing
# rows 198-221 down to rows 199-222) and populates the new row with the
# latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 198; Excel shifts rows 198:221 down to 199:222 and
# copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A198").Value = 11
$ws.Range("B198").Value = "Vega Monumental Concepción"
$ws.Range("C198").Value = "Bíobío"
$ws.Range("D198").Value = 44748
$ws.Range("E198").Value = 8
$ws.Range("F198").Value = 100112045
$ws.Range("G198").Value = "Zapallo"
$ws.Range("H198").Value = "Camote"
$ws.Range("I198").Value = "1a (guarda)"
$ws.Range("J198").Value = 550
$ws.Range("K198").Value = 500
$ws.Range("L198").Value = 600
$ws.Range("M198").Value = 545
$ws.Range("N198").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O198").Value = "Región de O'Higgins"
$ws.Range("P198").Value = 545
$ws.Range("Q198").Value = 1
$ws.Range("R198").Value = "Hortaliza"
